$wb = $excel.ActiveWorkbook

# --- gdp_calibrate sheet: add a new row (year 730, value 3000) ---
$gdp = $wb.Worksheets.Item("gdp_calibrate")
$gdp.Range("A7").Value = "Westeros"
$gdp.Range("B7").Value = 730
$gdp.Range("C7").Value = 3000
$gdp.Range("D7").Value = "T$"
$gdp.Range("E7").Value = "not sure what to put as variable GDP DNE for Westeros"

# --- config sheet: add a "year" column with desirable years ---
$config = $wb.Worksheets.Item("config")
$config.Range("E1").Value = "year"
$config.Range("E2").Value = 700
$config.Range("E3").Value = 710
$config.Range("E4").Value = 720
